$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.654872894287109
$ws.Range("B1").Value = 2.490117073059082
$ws.Range("C1").Value = 1.784607648849487
$ws.Range("D1").Value = 1.6317138671875
$ws.Range("E1").Value = 1.65909481048584
